$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Cont Type (*)"
$ws.Range("B1").Value = "Cont Qty (*)"
